$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2503.413669864819, 1825.8615987511057, 1805.435296458116),
    @(2283.7042985547573, 1493.8527053304535, 1710.1725806226718),
    @(2560.4095235270556, 1801.9440272403326, 1920.2388694916276),
    @(2537.9421661743072, 1799.864344431141, 1646.2595615316329),
    @(2681.6711769455369, 1745.3775111589102, 1711.6271243643137),
    @(2529.4621565535494, 1702.2630343354697, 2011.3320419227862),
    @(2396.4138229327441, 2114.9264906491089, 2077.3148737878214),
    @(2466.5617480975629, 2066.6417710190935, 2063.6874026853902),
    @(2646.5013805246285, 1942.8177475121645, 1637.0888692757562),
    @(2445.758014446571, 1435.6782207248878, 1576.9788510410804),
    @(2163.7686519464837, 1676.8451866028288, 1606.5592030197213),
    @(2859.3559346016937, 2218.6895762930831, 1916.9616095439039),
    @(2574.3872249153419, 2072.2152191465848, 2070.4027135578858),
    @(2609.5010456619448, 2265.6740662602192, 2132.4259061940784),
    @(2575.2390392833559, 1954.5493876737485, 2192.9488406502014),
    @(2591.9502308411556, 1797.3965370695403, 1572.5026384786988)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
